$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace employee ids with employee names (string lookups instead of numeric ids)
$ws.Range("A2").Value = "Elijah Kyule"
$ws.Range("A3").Value = "Christopher Ndemi"
$ws.Range("A4").Value = "Tobias Mwalili"

# Update payroll figures for row 4 (Tobias Mwalili)
$ws.Range("B4").Value = 200000
$ws.Range("C4").Value = 60000
$ws.Range("D4").Value = 9000

# Widen column A to fit the longer employee names
$ws.Columns.Item(1).ColumnWidth = 22.3

# Update the active cell selection
$ws.Range("E4").Select()
